$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "08cef136-e4b6-497d-9fac-d7ababb06d4d"
$ws.Range("C2").Value = "2024-07-19T07:43:00.000Z"
$ws.Range("D2").Value = "2024-07-19T08:00:00.000Z"
$ws.Range("I2").Value = "https://www.notion.so/1-08cef136e4b6497d9facd7ababb06d4d"
$ws.Range("L2").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B3").Value = "2c2ffa37-24de-460f-8876-8653224f2ff0"
$ws.Range("C3").Value = "2024-07-17T12:18:00.000Z"
$ws.Range("D3").Value = "2024-07-19T08:00:00.000Z"
$ws.Range("I3").Value = "https://www.notion.so/2-2c2ffa3724de460f88768653224f2ff0"
$ws.Range("L3").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B4").Value = "cced709d-cbb3-4950-963d-64748a3f96c9"
$ws.Range("C4").Value = "2024-07-17T12:18:00.000Z"
$ws.Range("D4").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I4").Value = "https://www.notion.so/3-cced709dcbb34950963d64748a3f96c9"
$ws.Range("L4").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B5").Value = "f1d8ccda-f637-4f2d-915b-92da1400d99d"
$ws.Range("C5").Value = "2024-07-17T12:18:00.000Z"
$ws.Range("D5").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I5").Value = "https://www.notion.so/4-f1d8ccdaf6374f2d915b92da1400d99d"
$ws.Range("L5").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B6").Value = "8768deab-b9a3-473d-ba15-f02758483fdb"
$ws.Range("C6").Value = "2024-07-13T12:44:00.000Z"
$ws.Range("D6").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I6").Value = "https://www.notion.so/5-8768deabb9a3473dba15f02758483fdb"
$ws.Range("L6").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B7").Value = "222150b9-53a1-49d6-b41d-de3a0437e170"
$ws.Range("C7").Value = "2024-07-12T11:25:00.000Z"
$ws.Range("D7").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I7").Value = "https://www.notion.so/6-222150b953a149d6b41dde3a0437e170"
$ws.Range("L7").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B8").Value = "2d8adca0-1726-4176-bfc3-e44a4fe97b13"
$ws.Range("C8").Value = "2024-07-11T14:41:00.000Z"
$ws.Range("D8").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I8").Value = "https://www.notion.so/7-2d8adca017264176bfc3e44a4fe97b13"
$ws.Range("L8").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B9").Value = "7e8943d7-33eb-44a0-96b4-c7370eb06873"
$ws.Range("C9").Value = "2024-07-08T15:14:00.000Z"
$ws.Range("D9").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I9").Value = "https://www.notion.so/8-7e8943d733eb44a096b4c7370eb06873"
$ws.Range("L9").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B10").Value = "8ea311fb-b86c-4d56-ad05-7fed94812a32"
$ws.Range("C10").Value = "2024-07-07T14:57:00.000Z"
$ws.Range("D10").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I10").Value = "https://www.notion.so/9-8ea311fbb86c4d56ad057fed94812a32"
$ws.Range("L10").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B11").Value = "1060fd8d-47b8-4f3a-8ec7-b5bfd89ec67e"
$ws.Range("C11").Value = "2024-07-06T16:05:00.000Z"
$ws.Range("D11").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I11").Value = "https://www.notion.so/10-1060fd8d47b84f3a8ec7b5bfd89ec67e"
$ws.Range("L11").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B12").Value = "3185c0ad-0c6d-45f3-9ab0-17634502cc55"
$ws.Range("C12").Value = "2024-07-06T15:49:00.000Z"
$ws.Range("D12").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I12").Value = "https://www.notion.so/11-3185c0ad0c6d45f39ab017634502cc55"
$ws.Range("L12").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B13").Value = "d54f1efa-a863-42d5-9189-765fd4983cec"
$ws.Range("C13").Value = "2024-07-03T13:26:00.000Z"
$ws.Range("D13").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I13").Value = "https://www.notion.so/12-d54f1efaa86342d59189765fd4983cec"
$ws.Range("L13").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B14").Value = "52f1c71c-fb9a-4a52-b420-d5ba78d84357"
$ws.Range("C14").Value = "2024-07-03T13:20:00.000Z"
$ws.Range("D14").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I14").Value = "https://www.notion.so/13-52f1c71cfb9a4a52b420d5ba78d84357"
$ws.Range("L14").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B15").Value = "7791c7c1-0d0e-45e8-8da4-065604bcd21f"
$ws.Range("C15").Value = "2024-07-01T04:19:00.000Z"
$ws.Range("D15").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I15").Value = "https://www.notion.so/14-7791c7c10d0e45e88da4065604bcd21f"
$ws.Range("L15").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B16").Value = "00222ff8-cefb-4ae8-9158-3d12bac0731a"
$ws.Range("C16").Value = "2024-06-30T15:24:00.000Z"
$ws.Range("D16").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I16").Value = "https://www.notion.so/15-00222ff8cefb4ae891583d12bac0731a"
$ws.Range("L16").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B17").Value = "6187aa46-68a2-4a2d-8045-4548c4b66af9"
$ws.Range("C17").Value = "2024-06-30T15:24:00.000Z"
$ws.Range("D17").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I17").Value = "https://www.notion.so/16-6187aa4668a24a2d80454548c4b66af9"
$ws.Range("L17").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B18").Value = "464a57ea-ba65-4c01-9124-e01717e3f099"
$ws.Range("C18").Value = "2024-06-30T15:24:00.000Z"
$ws.Range("D18").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I18").Value = "https://www.notion.so/17-464a57eaba654c019124e01717e3f099"
$ws.Range("L18").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B19").Value = "94d7f750-1494-430b-8178-f828dd9b7995"
$ws.Range("C19").Value = "2024-06-30T15:24:00.000Z"
$ws.Range("D19").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I19").Value = "https://www.notion.so/18-94d7f7501494430b8178f828dd9b7995"
$ws.Range("L19").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B20").Value = "594b617c-1c12-4b59-a161-3b855ace2697"
$ws.Range("C20").Value = "2024-06-30T15:24:00.000Z"
$ws.Range("D20").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I20").Value = "https://www.notion.so/19-594b617c1c124b59a1613b855ace2697"
$ws.Range("L20").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B21").Value = "2f779725-2184-4679-954f-d9057ac30b88"
$ws.Range("C21").Value = "2024-06-30T15:24:00.000Z"
$ws.Range("D21").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I21").Value = "https://www.notion.so/20-2f77972521844679954fd9057ac30b88"
$ws.Range("L21").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B22").Value = "d2550ac2-132f-49ff-ae25-eaca9136f3e4"
$ws.Range("C22").Value = "2024-06-26T11:55:00.000Z"
$ws.Range("D22").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I22").Value = "https://www.notion.so/21-d2550ac2132f49ffae25eaca9136f3e4"
$ws.Range("L22").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B23").Value = "004ae277-0f96-4f1c-acb0-13a1cf8d7cc0"
$ws.Range("C23").Value = "2024-06-26T01:42:00.000Z"
$ws.Range("D23").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I23").Value = "https://www.notion.so/22-004ae2770f964f1cacb013a1cf8d7cc0"
$ws.Range("L23").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B24").Value = "3c276d55-7b72-423e-8f84-5129c576bb8c"
$ws.Range("C24").Value = "2024-06-26T01:41:00.000Z"
$ws.Range("D24").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I24").Value = "https://www.notion.so/23-3c276d557b72423e8f845129c576bb8c"
$ws.Range("L24").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B25").Value = "6b50a32b-5db4-4167-be33-8a430747a091"
$ws.Range("C25").Value = "2024-06-26T01:41:00.000Z"
$ws.Range("D25").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I25").Value = "https://www.notion.so/24-6b50a32b5db44167be338a430747a091"
$ws.Range("L25").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B26").Value = "39facd95-bd4d-4426-90e7-8a4d9968d5e3"
$ws.Range("C26").Value = "2024-06-26T01:41:00.000Z"
$ws.Range("D26").Value = "2024-07-19T08:01:00.000Z"
$ws.Range("I26").Value = "https://www.notion.so/25-39facd95bd4d442690e78a4d9968d5e3"
$ws.Range("L26").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B27").Value = "584ac372-86fe-4fcb-85b9-7b563f2d0fbd"
$ws.Range("C27").Value = "2024-06-26T01:41:00.000Z"
$ws.Range("D27").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I27").Value = "https://www.notion.so/26-584ac37286fe4fcb85b97b563f2d0fbd"
$ws.Range("L27").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B28").Value = "9e960e79-60a7-46dd-9260-80dea0b94e8e"
$ws.Range("C28").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D28").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I28").Value = "https://www.notion.so/27-9e960e7960a746dd926080dea0b94e8e"
$ws.Range("L28").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B29").Value = "82ce8aa2-3f83-4a5a-8bb4-8b6e0fa431d4"
$ws.Range("C29").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D29").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I29").Value = "https://www.notion.so/28-82ce8aa23f834a5a8bb48b6e0fa431d4"
$ws.Range("L29").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B30").Value = "f73d7cb7-2537-4af8-97a2-02fc4ab15d9e"
$ws.Range("C30").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D30").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I30").Value = "https://www.notion.so/29-f73d7cb725374af897a202fc4ab15d9e"
$ws.Range("L30").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B31").Value = "96b2ae74-933b-4a59-b431-8f36dd90de80"
$ws.Range("C31").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D31").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I31").Value = "https://www.notion.so/30-96b2ae74933b4a59b4318f36dd90de80"
$ws.Range("L31").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B32").Value = "8ddf806c-8a8b-4528-b218-2039b0175052"
$ws.Range("C32").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D32").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I32").Value = "https://www.notion.so/31-8ddf806c8a8b4528b2182039b0175052"
$ws.Range("L32").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B33").Value = "2990b1ef-c2ec-4fc0-a44d-b8bac4db74b9"
$ws.Range("C33").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D33").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I33").Value = "https://www.notion.so/32-2990b1efc2ec4fc0a44db8bac4db74b9"
$ws.Range("L33").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B34").Value = "7f37cb5c-5821-42ea-a19b-8910fba84fb5"
$ws.Range("C34").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D34").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I34").Value = "https://www.notion.so/33-7f37cb5c582142eaa19b8910fba84fb5"
$ws.Range("L34").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B35").Value = "0b7379bc-f352-4302-a3e8-8df0ea8dc0d5"
$ws.Range("C35").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D35").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I35").Value = "https://www.notion.so/34-0b7379bcf3524302a3e88df0ea8dc0d5"
$ws.Range("L35").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B36").Value = "b379d1e4-550e-4cb2-af0c-b6810ff97bb8"
$ws.Range("C36").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D36").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I36").Value = "https://www.notion.so/35-b379d1e4550e4cb2af0cb6810ff97bb8"
$ws.Range("L36").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B37").Value = "3646bf89-7b7c-46b0-8c14-3ac95984359e"
$ws.Range("C37").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D37").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I37").Value = "https://www.notion.so/36-3646bf897b7c46b08c143ac95984359e"
$ws.Range("L37").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B38").Value = "b10005d8-518d-4edd-bc59-8449a5129fd8"
$ws.Range("C38").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D38").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I38").Value = "https://www.notion.so/37-b10005d8518d4eddbc598449a5129fd8"
$ws.Range("L38").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B39").Value = "f06f8dd5-e07d-4056-914b-088f2bf058fa"
$ws.Range("C39").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D39").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I39").Value = "https://www.notion.so/38-f06f8dd5e07d4056914b088f2bf058fa"
$ws.Range("L39").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B40").Value = "380732e8-c5bd-46b8-bff1-9fbd6f6c42f4"
$ws.Range("C40").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D40").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I40").Value = "https://www.notion.so/39-380732e8c5bd46b8bff19fbd6f6c42f4"
$ws.Range("L40").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B41").Value = "e33e5d89-8000-4fdf-81c1-18a3d5babcb5"
$ws.Range("C41").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D41").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I41").Value = "https://www.notion.so/40-e33e5d8980004fdf81c118a3d5babcb5"
$ws.Range("L41").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B42").Value = "a43add4c-e37b-465c-a76a-0282423a6997"
$ws.Range("C42").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D42").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I42").Value = "https://www.notion.so/41-a43add4ce37b465ca76a0282423a6997"
$ws.Range("L42").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B43").Value = "0f6131e3-cc67-4a48-8b6d-de4ef7d1b40f"
$ws.Range("C43").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D43").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I43").Value = "https://www.notion.so/42-0f6131e3cc674a488b6dde4ef7d1b40f"
$ws.Range("L43").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B44").Value = "80e96903-d54f-4d06-ad2f-a980d77dccb2"
$ws.Range("C44").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D44").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I44").Value = "https://www.notion.so/43-80e96903d54f4d06ad2fa980d77dccb2"
$ws.Range("L44").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B45").Value = "ec05a2ee-b7a6-4a55-bb1f-a891a55add0a"
$ws.Range("C45").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D45").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I45").Value = "https://www.notion.so/44-ec05a2eeb7a64a55bb1fa891a55add0a"
$ws.Range("L45").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B46").Value = "c3840748-0026-48d1-85e6-db7edb3c0504"
$ws.Range("C46").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D46").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I46").Value = "https://www.notion.so/45-c3840748002648d185e6db7edb3c0504"
$ws.Range("L46").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B47").Value = "e41d1598-1d6c-4a38-a5ed-13a3bf7d0bc7"
$ws.Range("C47").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D47").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I47").Value = "https://www.notion.so/46-e41d15981d6c4a38a5ed13a3bf7d0bc7"
$ws.Range("L47").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B48").Value = "f935ce96-d7be-4cac-881b-b804231d5eea"
$ws.Range("C48").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D48").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I48").Value = "https://www.notion.so/47-f935ce96d7be4cac881bb804231d5eea"
$ws.Range("L48").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B49").Value = "af426fdc-06e6-4499-90da-573a73641241"
$ws.Range("C49").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D49").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I49").Value = "https://www.notion.so/48-af426fdc06e6449990da573a73641241"
$ws.Range("L49").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B50").Value = "84b1b189-fbd5-4f8f-b3f7-593faaa1042b"
$ws.Range("C50").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D50").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I50").Value = "https://www.notion.so/49-84b1b189fbd54f8fb3f7593faaa1042b"
$ws.Range("L50").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B51").Value = "2c707af5-923c-4344-98d5-e3bb47c7e443"
$ws.Range("C51").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D51").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I51").Value = "https://www.notion.so/50-2c707af5923c434498d5e3bb47c7e443"
$ws.Range("L51").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B52").Value = "801e99d3-7aef-4e90-b819-9626450680b6"
$ws.Range("C52").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D52").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I52").Value = "https://www.notion.so/51-801e99d37aef4e90b8199626450680b6"
$ws.Range("L52").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B53").Value = "589c8f49-bcb1-4fe9-9280-01f50e7f65e9"
$ws.Range("C53").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D53").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I53").Value = "https://www.notion.so/52-589c8f49bcb14fe9928001f50e7f65e9"
$ws.Range("L53").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B54").Value = "0952e540-9579-4e94-896c-937858a124d1"
$ws.Range("C54").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D54").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I54").Value = "https://www.notion.so/53-0952e54095794e94896c937858a124d1"
$ws.Range("L54").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B55").Value = "f5529a3c-c2cd-4cc3-9daf-faee511ed521"
$ws.Range("C55").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D55").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I55").Value = "https://www.notion.so/54-f5529a3cc2cd4cc39daffaee511ed521"
$ws.Range("L55").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B56").Value = "f492049b-ec0e-4cc3-ba2a-ff05d9bc1b9b"
$ws.Range("C56").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D56").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I56").Value = "https://www.notion.so/55-f492049bec0e4cc3ba2aff05d9bc1b9b"
$ws.Range("L56").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B57").Value = "1aa2fc39-ea35-41ca-8f39-b083de502a63"
$ws.Range("C57").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D57").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I57").Value = "https://www.notion.so/56-1aa2fc39ea3541ca8f39b083de502a63"
$ws.Range("L57").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B58").Value = "cb6000f4-2feb-4a98-8144-4628f0f82d9e"
$ws.Range("C58").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D58").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I58").Value = "https://www.notion.so/57-cb6000f42feb4a9881444628f0f82d9e"
$ws.Range("L58").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B59").Value = "e541355e-46b0-4e13-add8-1a769ba877c5"
$ws.Range("C59").Value = "2024-06-25T15:46:00.000Z"
$ws.Range("D59").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I59").Value = "https://www.notion.so/58-e541355e46b04e13add81a769ba877c5"
$ws.Range("L59").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B60").Value = "c515f008-b17f-48c5-878e-502c037a0818"
$ws.Range("C60").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D60").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I60").Value = "https://www.notion.so/59-c515f008b17f48c5878e502c037a0818"
$ws.Range("L60").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B61").Value = "434e47e8-b704-495b-a187-ae5b8a77bd8d"
$ws.Range("C61").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D61").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I61").Value = "https://www.notion.so/60-434e47e8b704495ba187ae5b8a77bd8d"
$ws.Range("L61").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B62").Value = "d6e37602-1487-486e-beb8-6c14b5f3edca"
$ws.Range("C62").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D62").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I62").Value = "https://www.notion.so/61-d6e376021487486ebeb86c14b5f3edca"
$ws.Range("L62").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B63").Value = "7e00ea89-7a0f-4e22-a8dc-2776c56702cd"
$ws.Range("C63").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D63").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I63").Value = "https://www.notion.so/62-7e00ea897a0f4e22a8dc2776c56702cd"
$ws.Range("L63").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B64").Value = "7f0916e4-0d04-4a1d-ae32-ef60f001f862"
$ws.Range("C64").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D64").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I64").Value = "https://www.notion.so/63-7f0916e40d044a1dae32ef60f001f862"
$ws.Range("L64").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B65").Value = "168276c3-ef55-4ff9-9b8e-ea4af91d3dbe"
$ws.Range("C65").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D65").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I65").Value = "https://www.notion.so/64-168276c3ef554ff99b8eea4af91d3dbe"
$ws.Range("L65").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B66").Value = "d13bcc52-f036-4a4f-82a2-2a6e310eee1e"
$ws.Range("C66").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D66").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I66").Value = "https://www.notion.so/65-d13bcc52f0364a4f82a22a6e310eee1e"
$ws.Range("L66").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B67").Value = "109069f6-0fa0-4a94-a905-95a32e8bf30d"
$ws.Range("C67").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D67").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I67").Value = "https://www.notion.so/66-109069f60fa04a94a90595a32e8bf30d"
$ws.Range("L67").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B68").Value = "7cafd8bf-3912-4c28-aebb-2c01fe4a14d9"
$ws.Range("C68").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D68").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I68").Value = "https://www.notion.so/67-7cafd8bf39124c28aebb2c01fe4a14d9"
$ws.Range("L68").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B69").Value = "c5aa38aa-b7c8-4ffa-89ee-93c581d320cb"
$ws.Range("C69").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D69").Value = "2024-07-19T08:02:00.000Z"
$ws.Range("I69").Value = "https://www.notion.so/68-c5aa38aab7c84ffa89ee93c581d320cb"
$ws.Range("L69").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B70").Value = "6daa2647-cb00-496e-b719-cd733569f662"
$ws.Range("C70").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D70").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I70").Value = "https://www.notion.so/69-6daa2647cb00496eb719cd733569f662"
$ws.Range("L70").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B71").Value = "fc89e627-ac81-4d73-8eb7-6c98511632f3"
$ws.Range("C71").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D71").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I71").Value = "https://www.notion.so/70-fc89e627ac814d738eb76c98511632f3"
$ws.Range("L71").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B72").Value = "111ad337-40f3-4b55-86b8-31259b1570fe"
$ws.Range("C72").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D72").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I72").Value = "https://www.notion.so/71-111ad33740f34b5586b831259b1570fe"
$ws.Range("L72").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B73").Value = "5c0f41ae-978e-4fab-a3e4-8af8d19c6f73"
$ws.Range("C73").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D73").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I73").Value = "https://www.notion.so/72-5c0f41ae978e4faba3e48af8d19c6f73"
$ws.Range("L73").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B74").Value = "856a3cd5-f46a-4018-9149-626d50d539ae"
$ws.Range("C74").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D74").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I74").Value = "https://www.notion.so/73-856a3cd5f46a40189149626d50d539ae"
$ws.Range("L74").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B75").Value = "3ef7f7a7-7144-4602-b6ff-afe3149c76f3"
$ws.Range("C75").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D75").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I75").Value = "https://www.notion.so/74-3ef7f7a771444602b6ffafe3149c76f3"
$ws.Range("L75").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B76").Value = "e8285e8e-bcad-4bf7-9bda-2bc436cc7c77"
$ws.Range("C76").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D76").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I76").Value = "https://www.notion.so/75-e8285e8ebcad4bf79bda2bc436cc7c77"
$ws.Range("L76").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B77").Value = "3d9b798d-e8d7-47da-afb7-e456bc579964"
$ws.Range("C77").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D77").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I77").Value = "https://www.notion.so/76-3d9b798de8d747daafb7e456bc579964"
$ws.Range("L77").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B78").Value = "39475552-6400-483e-a749-c503f5bacb66"
$ws.Range("C78").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D78").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I78").Value = "https://www.notion.so/77-394755526400483ea749c503f5bacb66"
$ws.Range("L78").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B79").Value = "4fcf70d6-d25a-43d6-bf81-8042b1523735"
$ws.Range("C79").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D79").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I79").Value = "https://www.notion.so/78-4fcf70d6d25a43d6bf818042b1523735"
$ws.Range("L79").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B80").Value = "1744c014-cf66-471c-8b3e-38faec3a07e9"
$ws.Range("C80").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D80").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I80").Value = "https://www.notion.so/79-1744c014cf66471c8b3e38faec3a07e9"
$ws.Range("L80").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B81").Value = "300d58e1-86d2-4cf3-a1a0-bf0d53b82bcf"
$ws.Range("C81").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D81").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I81").Value = "https://www.notion.so/80-300d58e186d24cf3a1a0bf0d53b82bcf"
$ws.Range("L81").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B82").Value = "99c24a5a-826f-4151-a69f-885ecdd28821"
$ws.Range("C82").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D82").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I82").Value = "https://www.notion.so/81-99c24a5a826f4151a69f885ecdd28821"
$ws.Range("L82").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B83").Value = "9b80c98a-f74d-4a5a-85c1-9876e46b346c"
$ws.Range("C83").Value = "2024-06-25T15:45:00.000Z"
$ws.Range("D83").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I83").Value = "https://www.notion.so/82-9b80c98af74d4a5a85c19876e46b346c"
$ws.Range("L83").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B84").Value = "68fe2bc2-21c0-48cb-ae81-1d7f6dcc7932"
$ws.Range("C84").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D84").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I84").Value = "https://www.notion.so/83-68fe2bc221c048cbae811d7f6dcc7932"
$ws.Range("L84").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B85").Value = "94cda438-3192-4b5f-a31c-c84e0ff157cd"
$ws.Range("C85").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D85").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I85").Value = "https://www.notion.so/84-94cda43831924b5fa31cc84e0ff157cd"
$ws.Range("L85").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B86").Value = "ca53893d-bfb7-4584-8516-39aef063440b"
$ws.Range("C86").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D86").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I86").Value = "https://www.notion.so/85-ca53893dbfb74584851639aef063440b"
$ws.Range("L86").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B87").Value = "3a542931-2b6e-40c6-b1ba-6c4531fbffa4"
$ws.Range("C87").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D87").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I87").Value = "https://www.notion.so/86-3a5429312b6e40c6b1ba6c4531fbffa4"
$ws.Range("L87").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B88").Value = "4d69d396-cf6e-42ef-b972-4b7ca1bb1cec"
$ws.Range("C88").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D88").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I88").Value = "https://www.notion.so/87-4d69d396cf6e42efb9724b7ca1bb1cec"
$ws.Range("L88").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B89").Value = "6b989ee4-9b14-4c20-a499-65208ae18ba5"
$ws.Range("C89").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D89").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I89").Value = "https://www.notion.so/88-6b989ee49b144c20a49965208ae18ba5"
$ws.Range("L89").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B90").Value = "7eb7a2aa-f1bb-4421-8aa4-518ecfff625e"
$ws.Range("C90").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D90").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I90").Value = "https://www.notion.so/89-7eb7a2aaf1bb44218aa4518ecfff625e"
$ws.Range("L90").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B91").Value = "2d371737-ff99-4b65-9a22-cbd894b2ac1f"
$ws.Range("C91").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D91").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I91").Value = "https://www.notion.so/90-2d371737ff994b659a22cbd894b2ac1f"
$ws.Range("L91").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B92").Value = "c5384494-b98f-4e9f-b016-3548157c2178"
$ws.Range("C92").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D92").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I92").Value = "https://www.notion.so/91-c5384494b98f4e9fb0163548157c2178"
$ws.Range("L92").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B93").Value = "12c869d9-2f55-4bff-8409-2cfc8fc5d429"
$ws.Range("C93").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D93").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I93").Value = "https://www.notion.so/92-12c869d92f554bff84092cfc8fc5d429"
$ws.Range("L93").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B94").Value = "ac5a0e65-c26d-4fe6-a273-a921c749d513"
$ws.Range("C94").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D94").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I94").Value = "https://www.notion.so/93-ac5a0e65c26d4fe6a273a921c749d513"
$ws.Range("L94").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B95").Value = "c3c0c752-885b-4329-9404-a29b0a9efd01"
$ws.Range("C95").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D95").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I95").Value = "https://www.notion.so/94-c3c0c752885b43299404a29b0a9efd01"
$ws.Range("L95").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B96").Value = "66c02225-dbe2-4eac-8bf7-ecdbdda09d48"
$ws.Range("C96").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D96").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I96").Value = "https://www.notion.so/95-66c02225dbe24eac8bf7ecdbdda09d48"
$ws.Range("L96").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B97").Value = "94a638d5-726d-4740-9797-8e4f3775a742"
$ws.Range("C97").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D97").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I97").Value = "https://www.notion.so/96-94a638d5726d474097978e4f3775a742"
$ws.Range("L97").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B98").Value = "916f4c61-6a8f-4943-b6c0-1a7b46ee4cde"
$ws.Range("C98").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D98").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I98").Value = "https://www.notion.so/97-916f4c616a8f4943b6c01a7b46ee4cde"
$ws.Range("L98").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B99").Value = "dc936191-5fd9-42da-8682-caa4e1831593"
$ws.Range("C99").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D99").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I99").Value = "https://www.notion.so/98-dc9361915fd942da8682caa4e1831593"
$ws.Range("L99").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B100").Value = "98971b82-bc7e-406c-b5a8-32f078924e5b"
$ws.Range("C100").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D100").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I100").Value = "https://www.notion.so/99-98971b82bc7e406cb5a832f078924e5b"
$ws.Range("L100").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B101").Value = "1dae0078-b2a7-4998-aee4-1469b4c89311"
$ws.Range("C101").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D101").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I101").Value = "https://www.notion.so/100-1dae0078b2a74998aee41469b4c89311"
$ws.Range("L101").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B102").Value = "c22f73f9-4cfe-41f0-b1e8-d20983b5801c"
$ws.Range("C102").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D102").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I102").Value = "https://www.notion.so/101-c22f73f94cfe41f0b1e8d20983b5801c"
$ws.Range("L102").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B103").Value = "e1eea6d7-d5a7-40ae-b005-ace756e63003"
$ws.Range("C103").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D103").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I103").Value = "https://www.notion.so/102-e1eea6d7d5a740aeb005ace756e63003"
$ws.Range("L103").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B104").Value = "b12a2dc1-7743-4fe1-8a66-784604ceff83"
$ws.Range("C104").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D104").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I104").Value = "https://www.notion.so/103-b12a2dc177434fe18a66784604ceff83"
$ws.Range("L104").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B105").Value = "a896264d-d720-447c-9605-ce9aed400e41"
$ws.Range("C105").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D105").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I105").Value = "https://www.notion.so/104-a896264dd720447c9605ce9aed400e41"
$ws.Range("L105").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B106").Value = "cbb64990-6e53-4766-9c60-ebbb82343224"
$ws.Range("C106").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D106").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I106").Value = "https://www.notion.so/105-cbb649906e5347669c60ebbb82343224"
$ws.Range("L106").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B107").Value = "4ec5bc4b-bd02-4d8d-a967-12b00846611f"
$ws.Range("C107").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D107").Value = "2024-07-19T08:03:00.000Z"
$ws.Range("I107").Value = "https://www.notion.so/106-4ec5bc4bbd024d8da96712b00846611f"
$ws.Range("L107").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B108").Value = "1ca34285-7b7b-45b3-8d63-0db08e839147"
$ws.Range("C108").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D108").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I108").Value = "https://www.notion.so/107-1ca342857b7b45b38d630db08e839147"
$ws.Range("L108").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B109").Value = "726c6db9-1331-49a4-ac95-b8bb50a908cf"
$ws.Range("C109").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D109").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I109").Value = "https://www.notion.so/108-726c6db9133149a4ac95b8bb50a908cf"
$ws.Range("L109").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B110").Value = "995f8767-dacb-4c2d-b37f-4a15b392ea08"
$ws.Range("C110").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D110").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I110").Value = "https://www.notion.so/109-995f8767dacb4c2db37f4a15b392ea08"
$ws.Range("L110").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B111").Value = "54c06f86-85fe-495c-8eed-c14761bc69d2"
$ws.Range("C111").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D111").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I111").Value = "https://www.notion.so/110-54c06f8685fe495c8eedc14761bc69d2"
$ws.Range("L111").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B112").Value = "4cf0cc74-afa5-4872-bc1e-a4b64ec04118"
$ws.Range("C112").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D112").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I112").Value = "https://www.notion.so/111-4cf0cc74afa54872bc1ea4b64ec04118"
$ws.Range("L112").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B113").Value = "2b9d2b03-2235-456c-b073-7d2b8cefe8a1"
$ws.Range("C113").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D113").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I113").Value = "https://www.notion.so/112-2b9d2b032235456cb0737d2b8cefe8a1"
$ws.Range("L113").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B114").Value = "a0afcb05-dbae-4f43-9611-7c939a084d20"
$ws.Range("C114").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D114").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I114").Value = "https://www.notion.so/113-a0afcb05dbae4f4396117c939a084d20"
$ws.Range("L114").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B115").Value = "74b364b1-15ec-4b14-aff9-65d64ce12e12"
$ws.Range("C115").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D115").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I115").Value = "https://www.notion.so/114-74b364b115ec4b14aff965d64ce12e12"
$ws.Range("L115").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B116").Value = "3c225fea-a1de-4673-aed9-f2ec10bf6e31"
$ws.Range("C116").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D116").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I116").Value = "https://www.notion.so/115-3c225feaa1de4673aed9f2ec10bf6e31"
$ws.Range("L116").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B117").Value = "88832275-d56a-4d93-9d6b-2189a13085b3"
$ws.Range("C117").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D117").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I117").Value = "https://www.notion.so/116-88832275d56a4d939d6b2189a13085b3"
$ws.Range("L117").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B118").Value = "9d451a84-6b7c-4019-a31c-9fed17f25dec"
$ws.Range("C118").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D118").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I118").Value = "https://www.notion.so/117-9d451a846b7c4019a31c9fed17f25dec"
$ws.Range("L118").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B119").Value = "0b237756-0adf-4566-9ceb-29df4911a1fc"
$ws.Range("C119").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D119").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I119").Value = "https://www.notion.so/118-0b2377560adf45669ceb29df4911a1fc"
$ws.Range("L119").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B120").Value = "aef2fde4-29f5-41a1-a393-4191f3aab1fb"
$ws.Range("C120").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D120").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I120").Value = "https://www.notion.so/119-aef2fde429f541a1a3934191f3aab1fb"
$ws.Range("L120").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B121").Value = "a5805d2a-974e-494d-8547-db6202ff0e8b"
$ws.Range("C121").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D121").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I121").Value = "https://www.notion.so/120-a5805d2a974e494d8547db6202ff0e8b"
$ws.Range("L121").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B122").Value = "907bfe00-b515-4d8d-896f-259cdf7d2552"
$ws.Range("C122").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D122").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I122").Value = "https://www.notion.so/121-907bfe00b5154d8d896f259cdf7d2552"
$ws.Range("L122").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B123").Value = "251b433b-4baa-41e0-9ec5-737b077a8955"
$ws.Range("C123").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D123").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I123").Value = "https://www.notion.so/122-251b433b4baa41e09ec5737b077a8955"
$ws.Range("L123").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B124").Value = "5b406471-905f-4cb1-ab92-7b89d36b8114"
$ws.Range("C124").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D124").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I124").Value = "https://www.notion.so/123-5b406471905f4cb1ab927b89d36b8114"
$ws.Range("L124").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B125").Value = "aef49589-f864-4815-b926-fad034a6d746"
$ws.Range("C125").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D125").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I125").Value = "https://www.notion.so/124-aef49589f8644815b926fad034a6d746"
$ws.Range("L125").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B126").Value = "dfcce9f3-165b-4c27-9a80-34317c295020"
$ws.Range("C126").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D126").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I126").Value = "https://www.notion.so/125-dfcce9f3165b4c279a8034317c295020"
$ws.Range("L126").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B127").Value = "5b513ab5-ea12-45a2-b96b-4bf9df9995ed"
$ws.Range("C127").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D127").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I127").Value = "https://www.notion.so/126-5b513ab5ea1245a2b96b4bf9df9995ed"
$ws.Range("L127").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B128").Value = "1ca4b9be-2ea6-46e8-bee9-100028cec31e"
$ws.Range("C128").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D128").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I128").Value = "https://www.notion.so/127-1ca4b9be2ea646e8bee9100028cec31e"
$ws.Range("L128").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B129").Value = "a16599b2-66e4-45d6-ba0c-3d49dff6c687"
$ws.Range("C129").Value = "2024-06-25T15:44:00.000Z"
$ws.Range("D129").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I129").Value = "https://www.notion.so/128-a16599b266e445d6ba0c3d49dff6c687"
$ws.Range("L129").Value = "41cabcaf-915d-46a5-8eff-38727be27269"
$ws.Range("B130").Value = "01558575-577b-4a9b-9b04-922298f3d78c"
$ws.Range("C130").Value = "2024-06-25T08:25:00.000Z"
$ws.Range("D130").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I130").Value = "https://www.notion.so/129-01558575577b4a9b9b04922298f3d78c"
$ws.Range("L130").Value = "532a166e-c2d9-42ff-bed3-a363f43543fb"
$ws.Range("B131").Value = "063fd668-0835-49af-a935-70ff8074aa87"
$ws.Range("C131").Value = "2024-06-25T08:25:00.000Z"
$ws.Range("D131").Value = "2024-07-19T07:59:00.000Z"
$ws.Range("I131").Value = "https://www.notion.so/130-063fd668083549afa93570ff8074aa87"
$ws.Range("L131").Value = "532a166e-c2d9-42ff-bed3-a363f43543fb"
